$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" footer field (Slide Master + all
#    11 Custom Layouts) from 14/03/2024 -> 20/03/2024, as if the deck had been
#    reopened/saved in PowerPoint on the newer date.
# ---------------------------------------------------------------------------
function Update-DatePlaceholders($shapes, $newText) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.HasTextFrame -and $shp.PlaceholderFormat.Type -eq 16) {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes "20/03/2024"

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholders $layout.Shapes "20/03/2024"
}

# ---------------------------------------------------------------------------
# 2) Slide 9 ("Evaluation") - colour the bullet list in "Content Placeholder 2"
#    with the accent6 theme colour.
# ---------------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
for ($j = 1; $j -le $slide9.Shapes.Count; $j++) {
    $shp = $slide9.Shapes.Item($j)
    if ($shp.Name -eq "Content Placeholder 2") {
        $tr = $shp.TextFrame.TextRange
        $tr.Font.Color.ObjectThemeColor = 10  # msoThemeColorAccent6
    }
}
